$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "cryptos" price list: Price (D) and Volume(1h) (E) columns.
# Both columns hold text (the Price column keeps thousands separators and
# trailing zeros literally, so it must stay text, not be reinterpreted as a
# number); any new Price value that happens to look like a plain decimal is
# entered with a leading apostrophe so Excel keeps storing it as text
# instead of silently converting it (and dropping e.g. a trailing "0").

$ws.Range("D2").Value = "59.999.06"
$ws.Range("E2").Value = "  -0.95%  "
$ws.Range("D3").Value = "2.418.72"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'553.04"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +4.15%  "
$ws.Range("D9").Value = "'0.105"
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("D13").Value = "'25.33"
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").Value = "2.849.73"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").Value = "59.931.32"
$ws.Range("E15").Value = "  -0.87%  "
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("D17").Value = "2.407.56"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").Value = "'327.71"
$ws.Range("E20").Value = "  -2.39%  "
$ws.Range("E21").Value = "  -3.77%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").Value = "'66.03"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("E24").Value = "  +3.68%  "
$ws.Range("D25").Value = "'8.65"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("E26").Value = "  +0.17%  "
$ws.Range("D27").Value = "'1.40"
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("E28").Value = "  -1.99%  "
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("D30").Value = "'169.37"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("E31").Value = "  -4.09%  "
$ws.Range("E32").Value = "  +1.18%  "
$ws.Range("D33").Value = "'18.62"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E37").Value = "  -2.09%  "
$ws.Range("D38").Value = "'326.52"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("E39").Value = "  -2.39%  "
$ws.Range("D40").Value = "'0.406"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("E41").Value = "  -1.77%  "
$ws.Range("D42").Value = "'140.18"
$ws.Range("E42").Value = "  -2.97%  "
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("E44").Value = "  -0.88%  "
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").Value = "'0.402"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").Value = "'0.0223"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("E50").Value = "  -4.08%  "
$ws.Range("E51").Value = "  -0.99%  "
